$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 28724.416059566604
    3  = 96498.875125477236
    4  = 96484.654994095006
    5  = 81949.788452993496
    6  = 48389.937403419724
    7  = 104549.86787078655
    8  = 102228.83713563389
    9  = 81280.101526164945
    10 = 82870.820477338653
    11 = 87521.314502545225
    12 = 87480.365440534282
    13 = 121875.69746683838
    14 = 50698.799158333859
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 2).Value = $values[$row]
}
